$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 3007.889
$ws.Range("I4").Value = 3323.875
$ws.Range("K4").Value = 3323.875
$ws.Range("M4").Value = -3209.875

$ws.Range("H12").Value = 529.2
$ws.Range("J12").Value = 574.5
$ws.Range("L12").Value = 574.5
$ws.Range("N12").Value = -914.5

$ws.Range("H17").Value = 572.8246
$ws.Range("J17").Value = 572.8246
$ws.Range("L17").Value = 1718.4738
$ws.Range("N17").Value = -2054.4738

$ws.Range("H33").Value = 9091341
$ws.Range("I33").Value = 5882859
$ws.Range("K33").Value = 5882859
$ws.Range("M33").Value = -5882630

$ws.Range("H39").Value = 2625
$ws.Range("J39").Value = 2625
$ws.Range("L39").Value = 7875
$ws.Range("N39").Value = -8467

$ws.Range("H51").Value = 26079.7
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 26079.7
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 26079.7
$ws.Range("M51").Value = ""
$ws.Range("N51").Value = -27047.7

$ws.Range("H53").Value = 665.43335
$ws.Range("I53").Value = 69.92308
$ws.Range("J53").Value = 1120.8235
$ws.Range("K53").Value = 69.92308
$ws.Range("L53").Value = 1120.8235
$ws.Range("M53").Value = 567.07692
$ws.Range("N53").Value = -2394.8235

$ws.Range("H76").Value = 3850335
$ws.Range("I76").Value = 5886159.5
$ws.Range("K76").Value = 5886159.5
$ws.Range("M76").Value = -5885844.5

$ws.Range("H79").Value = 3850335
$ws.Range("I79").Value = 5886159.5
$ws.Range("K79").Value = 5886159.5
$ws.Range("M79").Value = -5885067.5

$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").Value = ""

$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").Value = ""

$ws.Range("H125").Value = 2699.9285
$ws.Range("I125").Value = 1820.5
$ws.Range("J125").Value = 3359.5
$ws.Range("K125").Value = 16384.5
$ws.Range("L125").Value = 30235.5
$ws.Range("M125").Value = -13924.5
$ws.Range("N125").Value = -35155.5

$ws.Range("H131").Value = 2389.111
$ws.Range("J131").Value = 650
$ws.Range("L131").Value = 1950
$ws.Range("N131").Value = -12030

$ws.Range("H132").Value = 1902.8718
$ws.Range("I132").Value = 1408.92
$ws.Range("J132").Value = 2784.9285
$ws.Range("K132").Value = 4226.76
$ws.Range("L132").Value = 8354.7855
$ws.Range("M132").Value = -1696.76
$ws.Range("N132").Value = -13414.7855

$ws.Range("H138").Value = 5193.1772
$ws.Range("J138").Value = 6324.3335
$ws.Range("L138").Value = 18973.0005
$ws.Range("N138").Value = -29253.0005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 200296
$ws.Range("I74").Value = 273322.66
$ws.Range("J74").Value = 7297
$ws.Range("K74").Value = 273322.66
$ws.Range("L74").Value = 7297
$ws.Range("M74").Value = -272448.66
$ws.Range("N74").Value = -9045

$ws.Range("H77").Value = 200296
$ws.Range("I77").Value = 273322.66
$ws.Range("J77").Value = 7297
$ws.Range("K77").Value = 1366613.3
$ws.Range("L77").Value = 36485
$ws.Range("M77").Value = -1362245.3
$ws.Range("N77").Value = -45221

$ws.Range("H97").Value = 1951777.5
$ws.Range("I97").Value = 2647412.2
$ws.Range("K97").Value = 2647412.2
$ws.Range("M97").Value = -2646916.2

$ws.Range("H102").Value = 1397.0476
$ws.Range("I102").Value = 1422.3684
$ws.Range("K102").Value = 1422.3684
$ws.Range("M102").Value = 199.6315999999999

$ws.Range("H110").Value = 30002338
$ws.Range("I110").Value = 52500884
$ws.Range("J110").Value = 4277.778
$ws.Range("K110").Value = 52500884
$ws.Range("L110").Value = 4277.778
$ws.Range("M110").Value = -52498839
$ws.Range("N110").Value = -8367.778

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3955.5576
$ws.Range("I134").Value = 1775.4054
$ws.Range("K134").Value = 5326.216200000001
$ws.Range("M134").Value = -2791.216200000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 17653.75
$ws.Range("I62").Value = 13495.143
$ws.Range("K62").Value = 13495.143
$ws.Range("M62").Value = -12871.143

$ws.Range("H65").Value = 17653.75
$ws.Range("I65").Value = 13495.143
$ws.Range("K65").Value = 67475.715
$ws.Range("M65").Value = -64355.715

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 869.4545000000001
$ws.Range("I14").Value = 869.4545000000001
$ws.Range("K14").Value = 2608.3635
$ws.Range("M14").Value = -2435.3635

$ws.Range("H47").Value = 143585.58
$ws.Range("I47").Value = 200519.8
$ws.Range("K47").Value = 601559.3999999999
$ws.Range("M47").Value = -601128.3999999999

$ws.Range("H48").Value = 3000
$ws.Range("J48").Value = 3000
$ws.Range("L48").Value = 9000
$ws.Range("N48").Value = -9500

$ws.Range("H68").Value = 51624.465
$ws.Range("J68").Value = 3157.2812
$ws.Range("L68").Value = 9471.8436
$ws.Range("N68").Value = -11093.8436

$ws.Range("H71").Value = 51624.465
$ws.Range("J71").Value = 3157.2812
$ws.Range("L71").Value = 28415.5308
$ws.Range("N71").Value = -36527.5308

$ws.Range("H131").Value = 13892923
$ws.Range("I131").Value = 23810196
$ws.Range("K131").Value = 71430588
$ws.Range("M131").Value = -71425548

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2426.2222
$ws.Range("I80").Value = 1991.3334
$ws.Range("J80").Value = 3296
$ws.Range("K80").Value = 1991.3334
$ws.Range("L80").Value = 3296
$ws.Range("M80").Value = -993.3334
$ws.Range("N80").Value = -5292

$ws.Range("H83").Value = 2426.2222
$ws.Range("I83").Value = 1991.3334
$ws.Range("J83").Value = 3296
$ws.Range("K83").Value = 9956.666999999999
$ws.Range("L83").Value = 16480
$ws.Range("M83").Value = -4964.666999999999
$ws.Range("N83").Value = -26464

$ws.Range("H132").Value = 6220.4707
$ws.Range("I132").Value = 3145.6667
$ws.Range("K132").Value = 9437.000100000001
$ws.Range("M132").Value = -6907.000100000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 61838.445
$ws.Range("I7").Value = 79776.30499999999
$ws.Range("K7").Value = 79776.30499999999
$ws.Range("M7").Value = -79664.30499999999

$ws.Range("H29").Value = 19997
$ws.Range("J29").Value = 19997
$ws.Range("L29").Value = 19997
$ws.Range("N29").Value = -20587

$ws.Range("H93").Value = 1389.2609
$ws.Range("J93").Value = 1961.7693
$ws.Range("L93").Value = 1961.7693
$ws.Range("N93").Value = -4457.7693

$ws.Range("H126").Value = 61838.445
$ws.Range("I126").Value = 79776.30499999999
$ws.Range("K126").Value = 239328.915
$ws.Range("M126").Value = -236858.915

$ws.Range("H132").Value = 4976.8
$ws.Range("I132").Value = 3972.28
$ws.Range("J132").Value = 9999.4
$ws.Range("K132").Value = 11916.84
$ws.Range("L132").Value = 29998.2
$ws.Range("M132").Value = -9386.84
$ws.Range("N132").Value = -35058.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 4450
$ws.Range("I32").Value = 4450
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 4450
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -4133
$ws.Range("N32").Value = ""

$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").Value = ""

$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").Value = ""

Write-Output "Applied all Zalera_Profits updates"